$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.498.06'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '1.914.10'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'245.17"
$ws.Range("E5").Value = '  +1.22%  '
$ws.Range("D6").Value = "'0.9991"
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").Value = "'0.4820"
$ws.Range("E7").Value = '  +2.39%  '
$ws.Range("D8").Value = "'0.2893"
$ws.Range("E8").Value = '  +1.22%  '
$ws.Range("D9").Value = "'0.06720"
$ws.Range("E9").Value = '  -1.37%  '
$ws.Range("D10").Value = "'111.41"
$ws.Range("E10").Value = '  +4.39%  '
$ws.Range("D11").Value = "'19.17"
$ws.Range("E11").Value = '  +4.59%  '
$ws.Range("D12").Value = '1.911.60'
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").Value = "'0.07555"
$ws.Range("E13").Value = '  -2.18%  '
$ws.Range("D14").Value = "'5.259"
$ws.Range("E14").Value = '  +0.63%  '
$ws.Range("D15").Value = "'0.6704"
$ws.Range("E15").Value = '  +1.76%  '
$ws.Range("D16").Value = "'286.88"
$ws.Range("E16").Value = '  -2.25%  '
$ws.Range("D17").Value = '30.489.14'
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").Value = "'0.000007614"
$ws.Range("E18").Value = '  -0.27%  '
$ws.Range("D19").Value = "'0.9993"
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").Value = "'12.89"
$ws.Range("E20").Value = '  -0.32%  '
$ws.Range("D21").Value = '2.162.35'
$ws.Range("E21").Value = '  +0.68%  '
$ws.Range("D22").Value = "'5.457"
$ws.Range("E22").Value = '  +4.60%  '
$ws.Range("D23").Value = "'0.9994"
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").Value = "'6.409"
$ws.Range("E24").Value = '  +3.17%  '
$ws.Range("D25").Value = "'9.449"
$ws.Range("E25").Value = '  +1.26%  '
$ws.Range("D26").Value = "'163.83"
$ws.Range("E26").Value = '  -2.57%  '
$ws.Range("D27").Value = "'20.29"
$ws.Range("E27").Value = '  -5.77%  '
$ws.Range("D28").Value = "'2.136"
$ws.Range("E28").Value = '  +2.78%  '
$ws.Range("E29").Value = '  -0.94%  '
$ws.Range("D30").Value = "'1.408"
$ws.Range("E30").Value = '  +2.89%  '
$ws.Range("D31").Value = "'4.161"
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("D32").Value = "'4.049"
$ws.Range("E32").Value = '  +1.46%  '
$ws.Range("E33").Value = '  -1.48%  '
$ws.Range("D34").Value = "'0.7275"
$ws.Range("E34").Value = '  -2.02%  '
$ws.Range("D35").Value = "'1.132"
$ws.Range("E35").Value = '  -1.85%  '
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").Value = "'0.02048"
$ws.Range("E37").Value = '  -1.99%  '
$ws.Range("D38").Value = "'2.722"
$ws.Range("E38").Value = '  -0.95%  '
$ws.Range("E39").Value = '  -0.40%  '
$ws.Range("D40").Value = "'110.67"
$ws.Range("E40").Value = '  +0.65%  '
$ws.Range("D41").Value = "'2.012"
$ws.Range("E41").Value = '  -2.78%  '
$ws.Range("D42").Value = "'0.4430"
$ws.Range("E42").Value = '  +3.56%  '
$ws.Range("D43").Value = "'0.8642"
$ws.Range("E43").Value = '  -0.71%  '
$ws.Range("D44").Value = "'5.878"
$ws.Range("E44").Value = '  +0.43%  '
$ws.Range("D45").Value = "'0.9993"
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").Value = "'67.96"
$ws.Range("E46").Value = '  +0.12%  '
$ws.Range("D47").Value = "'7.344"
$ws.Range("E47").Value = '  +2.28%  '
$ws.Range("D48").Value = "'48.67"
$ws.Range("E48").Value = '  -4.56%  '
$ws.Range("D49").Value = "'9.268"
$ws.Range("E49").Value = '  -0.78%  '
$ws.Range("D50").Value = "'0.1239"
$ws.Range("E50").Value = '  +2.18%  '
$ws.Range("D51").Value = "'34.81"
$ws.Range("E51").Value = '  -0.33%  '
